# gsc-export/HTTPS.xlsx -- roll the export window forward by one day.
#
# The source report always covers a trailing ~90-day window. On each
# refresh the oldest day drops off the front and a new day is appended
# at the end, and the "HTTPS URLs" series shifts up to match (the
# newest row repeats its previous count until the next refresh
# supplies a real figure for it).
#
# Old window: 2025-11-20 .. 2026-02-17 (rows 2..91)
# New window: 2025-11-21 .. 2026-02-18 (rows 2..91)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Boundary dates: $dates[0] is the old row-2 date, $dates[90] is the new
# row-91 date. Row r (2..91) moves from $dates[r-2] to $dates[r-1].
$dates = @("2025-11-20", "2025-11-21", "2025-11-22", "2025-11-23", "2025-11-24", "2025-11-25", "2025-11-26", "2025-11-27", "2025-11-28", "2025-11-29", "2025-11-30", "2025-12-01", "2025-12-02", "2025-12-03", "2025-12-04", "2025-12-05", "2025-12-06", "2025-12-07", "2025-12-08", "2025-12-09", "2025-12-10", "2025-12-11", "2025-12-12", "2025-12-13", "2025-12-14", "2025-12-15", "2025-12-16", "2025-12-17", "2025-12-18", "2025-12-19", "2025-12-20", "2025-12-21", "2025-12-22", "2025-12-23", "2025-12-24", "2025-12-25", "2025-12-26", "2025-12-27", "2025-12-28", "2025-12-29", "2025-12-30", "2025-12-31", "2026-01-01", "2026-01-02", "2026-01-03", "2026-01-04", "2026-01-05", "2026-01-06", "2026-01-07", "2026-01-08", "2026-01-09", "2026-01-10", "2026-01-11", "2026-01-12", "2026-01-13", "2026-01-14", "2026-01-15", "2026-01-16", "2026-01-17", "2026-01-18", "2026-01-19", "2026-01-20", "2026-01-21", "2026-01-22", "2026-01-23", "2026-01-24", "2026-01-25", "2026-01-26", "2026-01-27", "2026-01-28", "2026-01-29", "2026-01-30", "2026-01-31", "2026-02-01", "2026-02-02", "2026-02-03", "2026-02-04", "2026-02-05", "2026-02-06", "2026-02-07", "2026-02-08", "2026-02-09", "2026-02-10", "2026-02-11", "2026-02-12", "2026-02-13", "2026-02-14", "2026-02-15", "2026-02-16", "2026-02-17", "2026-02-18")

for ($row = 2; $row -le 91; $row++) {
    $newDate = $dates[$row - 1]
    $cell = $ws.Cells.Item($row, 1)
    # Leading "'" forces a literal text entry so a date-shaped string like
    # "2026-02-18" isn't reinterpreted as a date serial; ClearFormats()
    # then drops the quote-prefix formatting flag that introduces, so the
    # cell ends up back on the sheet's plain default style.
    $cell.Value = "'" + $newDate
    $cell.ClearFormats()
}

# "HTTPS URLs" counts for rows 2..91 -- the old C3..C91 values shifted up
# one row, with row 91 repeating its previous count as a placeholder.
$newCounts = @(25, 26, 26, 25, 25, 27, 28, 28, 27, 27, 27, 27, 27, 26, 25, 25, 25, 26, 27, 27, 29, 29, 30, 30, 31, 31, 31, 31, 31, 32, 32, 32, 32, 30, 31, 32, 30, 28, 28, 28, 28, 29, 29, 28, 27, 27, 28, 27, 27, 27, 27, 26, 26, 27, 26, 26, 25, 25, 25, 25, 26, 25, 24, 23, 24, 24, 24, 25, 26, 27, 28, 28, 28, 28, 28, 28, 28, 28, 27, 27, 28, 29, 29, 30, 30, 30, 31, 31, 30, 30)

for ($row = 2; $row -le 91; $row++) {
    $ws.Cells.Item($row, 3).Value = $newCounts[$row - 2]
}
